$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 131, shifting existing row 131 and below down by one
$ws.Rows.Item(131).EntireRow.Insert()

# Populate the new row 131 with data
$ws.Range("A131").Value = 11
$ws.Range("B131").Value = "Vega Monumental Concepción"
$ws.Range("C131").Value = "Bíobío"
$ws.Range("D131").Value = 44985
$ws.Range("E131").Value = 8
$ws.Range("F131").Value = 100112021
$ws.Range("G131").Value = "Ají"
$ws.Range("H131").Value = "Americana (o)"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 25
$ws.Range("K131").Value = 24000
$ws.Range("L131").Value = 25000
$ws.Range("M131").Value = 24600
$ws.Range("N131").Value = "$/saco 25 kilos"
$ws.Range("O131").Value = "Región Metropolitana"
$ws.Range("P131").Value = 984
$ws.Range("Q131").Value = 25
$ws.Range("R131").Value = "Hortaliza"
